$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (highlight fill) from the previous "Water" CF row (row 14)
# onto the two new rows before filling in their values.
$ws.Range("A14:C14").Copy() | Out-Null
$ws.Range("A15:C16").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row 15: Water / water::fossil well
$ws.Range("A15").Value = "Water"
$ws.Range("B15").Value = "water::fossil well"
$ws.Range("C15").Value = -0.000000001

# Row 16: Water / water::ground-, long-term
$ws.Range("A16").Value = "Water"
$ws.Range("B16").Value = "water::ground-, long-term"
$ws.Range("C16").Value = -0.000000001

# Match the resulting active selection recorded in the saved file.
$ws.Range("B19").Select() | Out-Null
